# Updates the Price (D) and Volume(1h) (E) columns of the cryptos list
# to the latest scraped values. Cells whose new Price text would be
# mis-parsed by Excel as a plain number (e.g. "0.999") are first marked
# as Text (NumberFormat "@") so they stay literal strings, matching the
# rest of the Price column (which already mixes "67.693.19"-style text
# that can't be auto-parsed as a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.693.19'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '3.286.82'
$ws.Range("E3").Value = '  -5.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.14'
$ws.Range("E5").Value = '  -3.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.74'
$ws.Range("E6").Value = '  -9.27%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '3.276.36'
$ws.Range("E8").Value = '  -5.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("E9").Value = '  -8.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("E10").Value = '  -11.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").Value = '  -2.34%  '
$ws.Range("E12").Value = '  -9.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.02'
$ws.Range("E13").Value = '  -12.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000248'
$ws.Range("E14").Value = '  -8.01%  '
$ws.Range("D15").Value = '3.810.75'
$ws.Range("E15").Value = '  -5.73%  '
$ws.Range("D16").Value = '67.688.87'
$ws.Range("E16").Value = '  -3.61%  '
$ws.Range("D17").Value = '3.284.97'
$ws.Range("E17").Value = '  -5.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '539.86'
$ws.Range("E18").Value = '  -8.78%  '
$ws.Range("E19").Value = '  -5.31%  '
$ws.Range("E20").Value = '  -12.18%  '
$ws.Range("E21").Value = '  -12.03%  '
$ws.Range("E22").Value = '  -11.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.93'
$ws.Range("E23").Value = '  -9.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.26'
$ws.Range("E24").Value = '  -10.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.71'
$ws.Range("E25").Value = '  -10.01%  '
$ws.Range("E26").Value = '  -9.97%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.21'
$ws.Range("E28").Value = '  -5.51%  '
$ws.Range("E29").Value = '  -12.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.54'
$ws.Range("E30").Value = '  -10.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.73'
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.17'
$ws.Range("E32").Value = '  -6.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.69'
$ws.Range("E33").Value = '  -15.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.85'
$ws.Range("E34").Value = '  -11.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '534.45'
$ws.Range("E35").Value = '  -6.52%  '
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  -6.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.64'
$ws.Range("E38").Value = '  -4.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0866'
$ws.Range("E39").Value = '  -10.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.10'
$ws.Range("E40").Value = '  -14.93%  '
$ws.Range("E41").Value = '  -9.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.86'
$ws.Range("E42").Value = '  -11.42%  '
$ws.Range("D43").Value = '2.966.01'
$ws.Range("E43").Value = '  -9.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.272'
$ws.Range("E44").Value = '  -9.51%  '
$ws.Range("D45").Value = '0.0₃0602'
$ws.Range("E45").Value = '  -14.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.23'
$ws.Range("E46").Value = '  -7.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.13'
$ws.Range("E47").Value = '  -12.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.38'
$ws.Range("E48").Value = '  -13.97%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("E50").Value = '  -9.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '123.98'
$ws.Range("E51").Value = '  -6.94%  '
